$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Hoja1")

# Step 1: bump the date in A1 by one day
$ws.Range("A1").Value = 45309

# Step 2: update unit prices in column D (rows 14-21)
$ws.Range("D14").Value = 76.307
$ws.Range("D15").Value = 99.298
$ws.Range("D16").Value = 118.919
$ws.Range("D17").Value = 190.27
$ws.Range("D18").Value = 237.838
$ws.Range("D19").Value = 297.297
$ws.Range("D20").Value = 330.991
$ws.Range("D21").Value = 366.666
